$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing data (D..K) to (F..M).
$ws.Range("D:E").Insert()

# The newly inserted D:E columns inherit column C's style; copy the number
# formats/styles from column F (the old column D, now shifted right) across
# the three contiguous data blocks so D:E match F (and thus G..M).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new quarters of data (columns D & E) plus the restated
# figure for the quarter that is now in column F (and, for a few rows whose
# source reports lagged further, columns G/H as well).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("D8").Value = 110200
$ws.Range("E8").Value = 108500
$ws.Range("F8").Value = 185200
$ws.Range("D9").Value = 32400
$ws.Range("E9").Value = 31000
$ws.Range("F9").Value = 51200
$ws.Range("D10").Value = 77800
$ws.Range("E10").Value = 77500
$ws.Range("F10").Value = 134000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 2600
$ws.Range("F14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("D17").Value = 46200
$ws.Range("E17").Value = 48900
$ws.Range("F17").Value = 78900
$ws.Range("D18").Value = 64000
$ws.Range("E18").Value = 59600
$ws.Range("F18").Value = 106300
$ws.Range("D20").Value = -11000
$ws.Range("E20").Value = 2700
$ws.Range("F20").Value = -1500
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("D23").Value = 53000
$ws.Range("E23").Value = 62200
$ws.Range("F23").Value = 104800
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("D26").Value = 53000
$ws.Range("E26").Value = 62200
$ws.Range("F26").Value = 104800
$ws.Range("D27").Value = 45000
$ws.Range("E27").Value = 54600
$ws.Range("F27").Value = 89600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("D32").Value = 11000
$ws.Range("E32").Value = -2700
$ws.Range("F32").Value = 1500
$ws.Range("D33").Value = 45000
$ws.Range("E33").Value = 54600
$ws.Range("F33").Value = 89600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("D35").Value = 45000
$ws.Range("E35").Value = 54600
$ws.Range("F35").Value = 89600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("D41").Value = 109800
$ws.Range("E41").Value = 99200
$ws.Range("F41").Value = 76400
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "NA"
$ws.Range("H43").Value = "NA"
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("D47").Value = 4928600
$ws.Range("E47").Value = 4828000
$ws.Range("F47").Value = 4866700
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "NA"
$ws.Range("H52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("D54").Value = 5095800
$ws.Range("E54").Value = 4974500
$ws.Range("F54").Value = 4981000
$ws.Range("D57").Value = 21500
$ws.Range("E57").Value = 4500
$ws.Range("F57").Value = 3500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("D59").Value = 93000
$ws.Range("E59").Value = 84600
$ws.Range("F59").Value = 87000
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("D61").Value = 2471500
$ws.Range("E61").Value = 2367000
$ws.Range("F61").Value = 2547900
$ws.Range("G61").Value = 1798700
$ws.Range("H61").Value = 1915700
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("D66").Value = 2586100
$ws.Range("E66").Value = 2456100
$ws.Range("F66").Value = 2638300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("D70").Value = 100
$ws.Range("E70").Value = 100
$ws.Range("F70").Value = 100
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("D72").Value = -130200
$ws.Range("E72").Value = -113600
$ws.Range("F72").Value = -106700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("D76").Value = 2509600
$ws.Range("E76").Value = 2518200
$ws.Range("F76").Value = 2342500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("D81").Value = 45000
$ws.Range("E81").Value = 54600
$ws.Range("F81").Value = 89600
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("F83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("D89").Value = 44800
$ws.Range("E89").Value = 49500
$ws.Range("F89").Value = 171600
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("F91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("D94").Value = -91500
$ws.Range("E94").Value = 45500
$ws.Range("F94").Value = -952900
$ws.Range("D96").Value = -50300
$ws.Range("E96").Value = -62300
$ws.Range("F96").Value = -114700
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("D100").Value = 57300
$ws.Range("E100").Value = -72200
$ws.Range("F100").Value = 780000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("D102").Value = 10600
$ws.Range("E102").Value = 22800
$ws.Range("F102").Value = -1300

Write-Output "Edit applied"
